$d = $word.ActiveDocument

# Right single quotation mark used in the Polish text ("branch'a", "branch'u").
$rsquo = [char]0x2019

$oldCaption = "Rys. 5. Dodanie pliku do branch" + $rsquo + "a."
$newCaption = "Rys. 6. Historia wersji pliku na branch" + $rsquo + "u."

# The last paragraph currently holds both the "Rys. 5" caption run and the
# _GoBack bookmark. Split it in two: a paragraph with just the caption text,
# and a paragraph with just the bookmark.
$found = $d.Content.Find.Execute($oldCaption, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $oldCaption + "^p", 2)

$bookmarkParaIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# Insert a blank paragraph between the caption paragraph and the bookmark
# paragraph.
$bookmarkPara.Range.InsertParagraphBefore()

# Re-fetch the bookmark paragraph (index shifted by one after the insert)
# and append a new paragraph after it for the "Rys. 6" caption.
$bookmarkParaIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

$tail = $bookmarkPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $newPara.Range
$insertPoint.Collapse(0)
$insertPoint.InsertBefore($newCaption)
